$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7207
$ws.Range("C3").Value = 179021
$ws.Range("C4").Value = 168977
$ws.Range("C8").Value = 64.97
